# Appending new "profile" test-case rows to the "Test Cases" sheet,
# matching the upstream commit "taking latest changes and appending new profile scripts".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# New test case rows (TCID, JIRA ID, Description, Runmode, Results)
$rows = @(
  @("TopicTypeaheadValidationTest", "TBD", "Verify that topic typeahead options should display while enter min 2 characters", "Y", "SKIP"),
  @("ProfileTabFocusTest", "TBD", "Verify that By default profile tab foucs should be on POST tab", "Y", "SKIP"),
  @("PostTabInfiniteScrollValidationTest", "TBD", "Verify that Post tab infinite scroll displaying the more available records", "Y", "SKIP"),
  @("CommentsTabInfiniteScrollValidationTest", "TBD", "Verify that Comments tab infinite scroll displaying the more available records", "Y", "PASS"),
  @("FollowersTabInfiniteScrollValidationTest", "TBD", "Verify that Followers tab infinite scroll displaying the more available records", "Y", "PASS"),
  @("FollowingTabInfiniteScrollValidationTest", "TBD", "Verify that Following tab infinite scroll displaying the more available records", "Y", "PASS")
)

$startRow = 30
for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $startRow + $i
  # Carry forward the row-29 formatting (borders/fonts/fills) onto the new row
  $ws.Range("A29:E29").Copy()
  $ws.Range("A" + $r + ":E" + $r).PasteSpecial(-4122)

  $data = $rows[$i]
  $ws.Range("A" + $r).Value = $data[0]
  $ws.Range("B" + $r).Value = $data[1]
  $ws.Range("C" + $r).Value = $data[2]
  $ws.Range("D" + $r).Value = $data[3]
  $ws.Range("E" + $r).Value = $data[4]
}

$excel.CutCopyMode = 0

$ws.Activate() | Out-Null
$ws.Range("D2:D35").Select() | Out-Null
